$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 126, shifting existing rows 126-156 down to 127-157.
$ws.Rows(126).Insert()

# Populate the newly inserted row 126 with its data.
$ws.Range("A126").Value = 9
$ws.Range("B126").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C126").Value = "Metropolitana"
$ws.Range("D126").Value = 44508
$ws.Range("E126").Value = 13
$ws.Range("F126").Value = 100112030
$ws.Range("G126").Value = "Poroto granado"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 16
$ws.Range("K126").Value = 34000
$ws.Range("L126").Value = 36000
$ws.Range("M126").Value = 35000
$ws.Range("N126").Value = "$/malla 25 kilos"
$ws.Range("O126").Value = "Región de Arica y Parinacota"
$ws.Range("P126").Value = 1400
$ws.Range("Q126").Value = 25
$ws.Range("R126").Value = "Hortaliza"
